# Add 2022-Q3 data
#
# 1) Insert a new quarterly sheet "2022-Q3" right after "总计" (i.e. before
#    the current "2022-Q2" sheet), by duplicating the "2022-Q2" sheet
#    (same column headers / styles) and then overwriting its data with the
#    2022-Q3 figures.
# 2) Insert a new row into the "总计" summary sheet for 2022-Q3, pushing the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet, positioned before "2022-Q2" ---
# Copy("$q2" as the "Before" sheet) inserts the duplicate immediately
# before "2022-Q2"; the worksheet reference then resolves to the newly
# inserted copy (same slot), so renaming it in place gives us "2022-Q3"
# right where we want it, with "2022-Q2" pushed one slot later.
$q2.Copy($q2)
$q3 = $q2
$q3.Name = "2022-Q3"

# Drop the extra data rows copied from 2022-Q2 (it had 10 data rows, Q3
# only has 4), keeping header row 1 and data rows 2-5.
$q3.Range("A6:H11").Delete()

# Columns B-G (fund code / name / size / position figures) are stored as
# text in the source data (e.g. fund code "003318" has a significant
# leading zero, and sizes like "10.25" are text, not numbers) - force the
# cells to Text format first so the COM layer doesn't silently reinterpret
# these as numeric values (which would also destroy the leading zeros).
$q3.Range("B2:G5").NumberFormat = "@"

# Overwrite the fund-holding data with the 2022-Q3 figures.
$q3.Range("B2").Value = "003318"
$q3.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$q3.Range("D2").Value = "10.25"
$q3.Range("E2").Value = "93.67"
$q3.Range("F2").Value = "1.25"
$q3.Range("G2").Value = "0.1281"
$q3.Range("H2").Value = 2

$q3.Range("B3").Value = "014133"
$q3.Range("C3").Value = "工银中证500六个月持有指数增强A"
$q3.Range("D3").Value = "1.64"
$q3.Range("E3").Value = "93.75"
$q3.Range("F3").Value = "1.64"
$q3.Range("G3").Value = "0.0269"
$q3.Range("H3").Value = 9

$q3.Range("B4").Value = "014134"
$q3.Range("C4").Value = "工银中证500六个月持有指数增强C"
$q3.Range("D4").Value = "0.90"
$q3.Range("E4").Value = "93.75"
$q3.Range("F4").Value = "1.64"
$q3.Range("G4").Value = "0.0148"
$q3.Range("H4").Value = 9

$q3.Range("B5").Value = "512260"
$q3.Range("C5").Value = "华安中证500行业中性低波动ETF"
$q3.Range("D5").Value = "1.07"
$q3.Range("E5").Value = "97.91"
$q3.Range("F5").Value = "1.31"
$q3.Range("G5").Value = "0.0140"
$q3.Range("H5").Value = 2

# --- 2. Insert the 2022-Q3 row into the "总计" summary sheet ---
# Range("X").Value reads back as an (unusable) property descriptor in this
# runtime rather than the live cell content, and Rows.Insert() fabricates a
# brand-new style index instead of reusing the sheet's existing ones — so
# rather than reading-and-shifting existing rows, or inserting a blank row,
# duplicate the row below (copies values *and* formatting in one shot) to
# make room, then overwrite every row with its final, literal value so the
# whole table ends up with one extra (2022-Q3) row and everything else
# pushed down by one, exactly as in the target data.
$total.Range("A7:D7").Copy($total.Range("A8:D8"))

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 2
$total.Range("D8").Value = 0.08

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 5
$total.Range("D7").Value = 0.19

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 0.13

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 0.44

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.79

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 10
$total.Range("D3").Value = 1.11

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.18
